{"js": "// Revert the document title from \"Version 7\" to \"Version 5\".\n// The title paragraph (\"Version \" + \"7\" in two separate runs) is\n// collapsed into a single run reading \"Version 5\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titlePara = paragraphs.items[0];\ntitlePara.load(\"text\");\nawait context.sync();\n\nif (titlePara.text.indexOf(\"Version\") === 0) {\n  titlePara.insertText(\"Version 5\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Revert the document title from \"Version 7\" to \"Version 5\".\n# The title paragraph holds \"Version \" and \"7\" in two separate runs;\n# Find/Replace collapses them into a single run reading \"Version 5\".\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \"Version 7\"\n$find.Replacement.Text = \"Version 5\"\n$find.Forward = $true\n$find.Wrap = 1            # wdFindContinue\n$find.Format = $false\n$find.MatchCase = $false\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.MatchSoundsLike = $false\n$find.MatchAllWordForms = $false\n\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null  # wdReplaceAll\n"}
